$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boss")

# Boss cycle data: name, HP threshold, multiplier (column C), repeated across 10 "days"
$bossData = @(
    ,@("飞龙", 6000000, 1)
    ,@("狂暴格里芬", 8000000, 1)
    ,@("兽人酋长", 10000000, 1.1)
    ,@("圣灵角鹿", 12000000, 1.1)
    ,@("牛头怪", 20000000, 1.2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
    ,@("飞龙", 6000000, 1.2)
    ,@("狂暴格里芬", 8000000, 1.2)
    ,@("兽人酋长", 10000000, 1.5)
    ,@("圣灵角鹿", 12000000, 1.7)
    ,@("牛头怪", 20000000, 2)
)

$r = 1
foreach ($row in $bossData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Activate()
[void]$ws.Range("K13").Select()
